$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column I
$ws.Range("I1").Value2 = "Other found locations"

# Update column E values (author lists) to the new variant (extra whitespace)
$ws.Range("E2").Value2 = "[Anthony V%Das%NULL%1,    Padmaja K%Rani%NULL%1,    Pravin K%Vaddavalli%NULL%1]"
$ws.Range("E3").Value2 = "[Gagan%Kalra%NULL%1,    Andrew M.%Williams%NULL%1,    Patrick W.%Commiskey%NULL%1,    Eve M. R.%Bowers%NULL%1,    Tadhg%Schempf%NULL%1,    José-Alain%Sahel%NULL%1,    Evan L.%Waxman%waxmane@upmc.edu%1,    Roxana%Fu%fur3@upmc.edu%1]"

# Add new column I values for each data row
$ws.Range("I2").Value2 = "_PMC"
$ws.Range("I3").Value2 = "_PMC_Springer"
